# Applies the cryptos.xlsx price/volume update described by the commit
# "Updated cryptos list ... with GitHub Actions".
# All target cells are plain text (inline/shared strings) in the workbook,
# so numeric-looking Price values are written with a leading apostrophe
# (quote-prefix) to force Excel to keep them as text instead of numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''61.363.55'
$ws.Range('E2').Value = '  -2.57%  '

# Row 3
$ws.Range('D3').Value = '''2.971.25'
$ws.Range('E3').Value = '  -2.46%  '

# Row 4
$ws.Range('E4').Value = '  +0.07%  '

# Row 5
$ws.Range('D5').Value = '''589.50'
$ws.Range('E5').Value = '  +1.47%  '

# Row 6
$ws.Range('D6').Value = '''141.96'
$ws.Range('E6').Value = '  -5.77%  '

# Row 7
$ws.Range('E7').Value = '  -0.01%  '

# Row 8
$ws.Range('E8').Value = '  -2.19%  '

# Row 9
$ws.Range('D9').Value = '''2.972.35'
$ws.Range('E9').Value = '  -2.56%  '

# Row 10
$ws.Range('E10').Value = '  -5.95%  '

# Row 11
$ws.Range('D11').Value = '''5.77'
$ws.Range('E11').Value = '  -0.45%  '

# Row 12
$ws.Range('E12').Value = '  +2.07%  '

# Row 13
$ws.Range('E13').Value = '  -3.01%  '

# Row 14
$ws.Range('D14').Value = '''33.97'
$ws.Range('E14').Value = '  -5.21%  '

# Row 15
$ws.Range('E15').Value = '  +1.34%  '

# Row 16
$ws.Range('D16').Value = '''3.463.85'
$ws.Range('E16').Value = '  -2.33%  '

# Row 17
$ws.Range('D17').Value = '''7.00'
$ws.Range('E17').Value = '  -1.22%  '

# Row 18
$ws.Range('D18').Value = '''61.310.30'
$ws.Range('E18').Value = '  -2.48%  '

# Row 19
$ws.Range('D19').Value = '''2.969.87'
$ws.Range('E19').Value = '  -2.33%  '

# Row 20
$ws.Range('D20').Value = '''447.86'
$ws.Range('E20').Value = '  -6.00%  '

# Row 21
$ws.Range('D21').Value = '''13.89'
$ws.Range('E21').Value = '  -2.26%  '

# Row 22
$ws.Range('D22').Value = '''0.683'
$ws.Range('E22').Value = '  -2.73%  '

# Row 23
$ws.Range('D23').Value = '''7.32'
$ws.Range('E23').Value = '  -2.10%  '

# Row 24
$ws.Range('D24').Value = '''81.15'
$ws.Range('E24').Value = '  -0.10%  '

# Row 25
$ws.Range('D25').Value = '''12.10'
$ws.Range('E25').Value = '  -3.63%  '

# Row 26
$ws.Range('E26').Value = '  -8.35%  '

# Row 27
$ws.Range('E27').Value = '  -0.01%  '

# Row 28
$ws.Range('D28').Value = '''9.93'
$ws.Range('E28').Value = '  -5.45%  '

# Row 29
$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  +0.03%  '

# Row 30
$ws.Range('E30').Value = '  +0.08%  '

# Row 31
$ws.Range('E31').Value = '  -5.88%  '

# Row 32
$ws.Range('E32').Value = '  -6.07%  '

# Row 33
$ws.Range('D33').Value = '''26.88'
$ws.Range('E33').Value = '  -2.60%  '

# Row 34
$ws.Range('E34').Value = '  -3.30%  '

# Row 35
$ws.Range('E35').Value = '  -4.50%  '

# Row 36
$ws.Range('D36').Value = '''0.0₃0775'
$ws.Range('E36').Value = '  -3.61%  '

# Row 37
$ws.Range('D37').Value = '''5.70'
$ws.Range('E37').Value = '  -2.63%  '

# Row 38
$ws.Range('D38').Value = '''50.12'
$ws.Range('E38').Value = '  -0.11%  '

# Row 39
$ws.Range('E39').Value = '  -5.38%  '

# Row 40
$ws.Range('D40').Value = '''9.11'
$ws.Range('E40').Value = '  +0.19%  '

# Row 41
$ws.Range('E41').Value = '  +4.80%  '

# Row 42
$ws.Range('D42').Value = '''2.75'
$ws.Range('E42').Value = '  -10.30%  '

# Row 43
$ws.Range('D43').Value = '''386.90'
$ws.Range('E43').Value = '  -8.73%  '

# Row 44
$ws.Range('D44').Value = '''0.0353'
$ws.Range('E44').Value = '  -1.69%  '

# Row 45
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = '''0.263'
$ws.Range('E45').Value = '  -7.36%  '

# Row 46
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '''2.696.80'
$ws.Range('E46').Value = '  -4.60%  '

# Row 47
$ws.Range('D47').Value = '''37.06'

# Row 48
$ws.Range('D48').Value = '''130.53'
$ws.Range('E48').Value = '  +2.22%  '

# Row 50
$ws.Range('D50').Value = '''0.108'
$ws.Range('E50').Value = '  -1.32%  '

# Row 51
$ws.Range('E51').Value = '  -0.89%  '
